$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.546.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.150.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.146.53"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.525"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  +0.43%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.663.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.118"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.368.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.137.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.716"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.25%  "
$ws.Range("E29").Value = "  +6.79%  "
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("E31").Value = "  -6.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.86%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -3.23%  "
$ws.Range("E35").Value = "  -1.33%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  +4.32%  "
$ws.Range("E38").Value = "  -1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.04"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "445.10"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0395"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.94%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.873.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.263"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "119.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.16%  "
